$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns Q (17) and R (18)
$ws.Range("Q1").Value = "PD"
$ws.Range("R1").Value = "N2"

# Copy the header formatting (bold, centered, bordered) from P1 onto Q1:R1
# so the new header cells share the same style as the rest of row 1.
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

# Per-row PD / N2 values for rows 2-45
$data = @(
    @('AGU-J', 'Fuera de Poligono OVL'),
    @('PUE-G', 'Fuera de Poligono OVL'),
    @('VCR-E', 'Fuera de Poligono OVL'),
    @('CLI-N', 'Fuera de Poligono OVL'),
    @('CON-H', 'Fuera de Poligono OVL'),
    @('RET-C', 'Fuera de Poligono OVL'),
    @('VCR-N', 'Fuera de Poligono OVL'),
    @('ALM-B', 'Fuera de Poligono OVL'),
    @('PCH-G', 'Fuera de Poligono OVL'),
    @('ATH-P', 'Fuera de Poligono OVL'),
    @('PUE-B', 'Fuera de Poligono OVL'),
    @('PPT-F', 'ARATO-25058.PO.1PPT'),
    @('PPT-K', 'Fuera de Poligono OVL'),
    @('ALM-G', 'Fuera de Poligono OVL'),
    @('ALM-F', 'Fuera de Poligono OVL'),
    @('CON-A', 'Fuera de Poligono OVL'),
    @('AGU-F', 'Fuera de Poligono OVL'),
    @('VCR-H', 'Fuera de Poligono OVL'),
    @('PAV-Q', 'Fuera de Poligono OVL'),
    @('COG-C', 'Fuera de Poligono OVL'),
    @('ALM-L', 'Fuera de Poligono OVL'),
    @('ALM-I', 'Fuera de Poligono OVL'),
    @('CON-B', 'Fuera de Poligono OVL'),
    @('CLI-L', 'Fuera de Poligono OVL'),
    @('PAV-T', 'Fuera de Poligono OVL'),
    @('CLI-B', 'Fuera de Poligono OVL'),
    @('ALM-J', 'Fuera de Poligono OVL'),
    @('ALM-C', 'Fuera de Poligono OVL'),
    @('PPT-S', 'Fuera de Poligono OVL'),
    @('ATH-J', 'Fuera de Poligono OVL'),
    @('NRA-A', 'Fuera de Poligono OVL'),
    @('AGU-B', 'Fuera de Poligono OVL'),
    @('PUE-O', 'ARATO-25058.PO.1PUE'),
    @('ATH-C', 'Fuera de Poligono OVL'),
    @('PUE-E', 'Fuera de Poligono OVL'),
    @('PUE-K', 'Fuera de Poligono OVL'),
    @('PAV-U', 'Fuera de Poligono OVL'),
    @('VCR-I', 'Fuera de Poligono OVL'),
    @('AGU-N', 'Fuera de Poligono OVL'),
    @('CON-E', 'Fuera de Poligono OVL'),
    @('BLO-G', 'Fuera de Poligono OVL'),
    @('RET-S', 'Fuera de Poligono OVL'),
    @('PUE-O', 'ARATO-25058.PO.1PUE'),
    @('BLO-B', 'Fuera de Poligono OVL')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $data[$i][0]
    $ws.Cells.Item($row, 18).Value = $data[$i][1]
}
